$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Rspo3"
$row2[0,2] = "Sdc4"
$row2[0,3] = "ECs"
$row2[0,4] = 1
$row2[0,5] = 0.3333333333333333
$row2[0,6] = 0.06729733333333333
$row2[0,7] = 0.201892
$row2[0,8] = 0.01373511018321553
$row2[0,9] = 0.01373511018321553
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 1.442875
$row2[0,13] = 4.328625
$row2[0,14] = 0.02047893724893121
$row2[0,15] = 0.02047893724893121
$row2[0,16] = 0.09710163983333332
$row2[0,17] = 0.8739147584999999
$row2[0,18] = 0.0002812804595492269
$row2[0,19] = 0.0002812804595492269
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Rspo3"
$row3[0,2] = "Sdc4"
$row3[0,3] = "FAPs"
$row3[0,4] = 1
$row3[0,5] = 0.3333333333333333
$row3[0,6] = 0.06729733333333333
$row3[0,7] = 0.201892
$row3[0,8] = 0.01373511018321553
$row3[0,9] = 0.01373511018321553
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 10.383857
$row3[0,13] = 31.151571
$row3[0,14] = 0.1473796107804731
$row3[0,15] = 0.1473796107804731
$row3[0,16] = 0.6988058858146667
$row3[0,17] = 6.289252972332
$row3[0,18] = 0.002024275192829218
$row3[0,19] = 0.002024275192829217
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Rspo3"
$row4[0,2] = "Sdc4"
$row4[0,3] = "MuSCs"
$row4[0,4] = 1
$row4[0,5] = 0.3333333333333333
$row4[0,6] = 0.06729733333333333
$row4[0,7] = 0.201892
$row4[0,8] = 0.01373511018321553
$row4[0,9] = 0.01373511018321553
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 27.934719
$row4[0,13] = 83.804157
$row4[0,14] = 0.3964815784233052
$row4[0,15] = 0.3964815784233051
$row4[0,16] = 1.879932096116
$row4[0,17] = 16.919388865044
$row4[0,18] = 0.005445718165259307
$row4[0,19] = 0.005445718165259305
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Rspo3"
$row5[0,2] = "Sdc4"
$row5[0,3] = "Resolving-Mac"
$row5[0,4] = 1
$row5[0,5] = 0.3333333333333333
$row5[0,6] = 0.06729733333333333
$row5[0,7] = 0.201892
$row5[0,8] = 0.01373511018321553
$row5[0,9] = 0.01373511018321553
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 30.695086
$row5[0,13] = 92.085258
$row5[0,14] = 0.4356598735472906
$row5[0,15] = 0.4356598735472905
$row5[0,16] = 2.065697434237333
$row5[0,17] = 18.591276908136
$row5[0,18] = 0.005983836365577783
$row5[0,19] = 0.005983836365577781
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Rspo3"
$row6[0,2] = "Sdc4"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 4.83236
$row6[0,7] = 14.49708
$row6[0,8] = 0.9862648898167845
$row6[0,9] = 0.9862648898167844
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 1.442875
$row6[0,13] = 4.328625
$row6[0,14] = 0.02047893724893121
$row6[0,15] = 0.02047893724893121
$row6[0,16] = 6.972491435
$row6[0,17] = 62.752422915
$row6[0,18] = 0.02019765678938198
$row6[0,19] = 0.02019765678938198
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Rspo3"
$row7[0,2] = "Sdc4"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 4.83236
$row7[0,7] = 14.49708
$row7[0,8] = 0.9862648898167845
$row7[0,9] = 0.9862648898167844
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 10.383857
$row7[0,13] = 31.151571
$row7[0,14] = 0.1473796107804731
$row7[0,15] = 0.1473796107804731
$row7[0,16] = 50.17853521252001
$row7[0,17] = 451.60681691268
$row7[0,18] = 0.1453553355876439
$row7[0,19] = 0.1453553355876439
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Rspo3"
$row8[0,2] = "Sdc4"
$row8[0,3] = "MuSCs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 4.83236
$row8[0,7] = 14.49708
$row8[0,8] = 0.9862648898167845
$row8[0,9] = 0.9862648898167844
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 27.934719
$row8[0,13] = 83.804157
$row8[0,14] = 0.3964815784233052
$row8[0,15] = 0.3964815784233051
$row8[0,16] = 134.99061870684
$row8[0,17] = 1214.91556836156
$row8[0,18] = 0.3910358602580459
$row8[0,19] = 0.3910358602580458
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Rspo3"
$row9[0,2] = "Sdc4"
$row9[0,3] = "Resolving-Mac"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 4.83236
$row9[0,7] = 14.49708
$row9[0,8] = 0.9862648898167845
$row9[0,9] = 0.9862648898167844
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 30.695086
$row9[0,13] = 92.085258
$row9[0,14] = 0.4356598735472906
$row9[0,15] = 0.4356598735472905
$row9[0,16] = 148.32970578296
$row9[0,17] = 1334.96735204664
$row9[0,18] = 0.4296760371817128
$row9[0,19] = 0.4296760371817127
$ws.Range("A9:T9").Value = $row9

Write-Host "Update complete"